# Applies the edits described by the commit:
# "Added support for longer quotes, fixed surplus numnber"
#
# 1. Fix the surplus multiplier in column K (rows 16, 18, 21, 24, 25, 29)
#    from 1.0565 down to 1.
# 2. Update the active cell selection on the sheet (from A29 to G5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Fix surplus number: change the K-column multiplier from 1.0565 to 1
# for each of the affected rows.
$ws.Range("K16").Value = 1
$ws.Range("K18").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("K24").Value = 1
$ws.Range("K25").Value = 1
$ws.Range("K29").Value = 1

# Move/update the current selection to G5 (reflecting the new active cell
# after the edits, e.g. scrolled up to support longer quotes).
$ws.Range("G5").Select()
